# Fruta / hortaliza, semanal
# Apply the weekly data refresh to the Damasco (Arica) price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 44533
$ws.Range("M2").Value = 140
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("Q2").Value = "$/caja 10 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1450
$ws.Range("T2").Value = 10

# --- Row 4 ---
$ws.Range("D4").Value = 44545
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 1361

# --- Row 5 ---
$ws.Range("D5").Value = 44174
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19500
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 1083

# --- Row 6 ---
$ws.Range("D6").Value = 44160
$ws.Range("L6").Value = "Primera"
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("S6").Value = 1361

# --- Row 7 ---
$ws.Range("D7").Value = 44544
$ws.Range("L7").Value = "Segunda"
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21000
$ws.Range("S7").Value = 1167

# --- Row 8 ---
$ws.Range("D8").Value = 44169
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 21000
$ws.Range("Q8").Value = "$/bandeja 18 kilos"
$ws.Range("R8").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S8").Value = 1167
$ws.Range("T8").Value = 18

Write-Host "Damasco (Arica) weekly refresh applied."
